$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 817.9107
$ws.Range("I15").Value = 817.9107
$ws.Range("K15").Value = 2453.7321
$ws.Range("M15").Value = -2284.7321

$ws.Range("H32").Value = 3664.7778
$ws.Range("I32").Value = 2913.878
$ws.Range("J32").Value = 6033
$ws.Range("K32").Value = 2913.878
$ws.Range("L32").Value = 6033
$ws.Range("M32").Value = -2587.878
$ws.Range("N32").Value = -6685

$ws.Range("H33").Value = 2112.5
$ws.Range("J33").Value = 873
$ws.Range("L33").Value = 873
$ws.Range("N33").Value = -1331

$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51248

$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -156240

$ws.Range("H69").Value = 362618.94
$ws.Range("I69").Value = 631765.06
$ws.Range("J69").Value = 3757.4167
$ws.Range("K69").Value = 1895295.18
$ws.Range("L69").Value = 11272.2501
$ws.Range("M69").Value = -1894421.18
$ws.Range("N69").Value = -13020.2501

$ws.Range("H72").Value = 362618.94
$ws.Range("I72").Value = 631765.06
$ws.Range("J72").Value = 3757.4167
$ws.Range("K72").Value = 5685885.540000001
$ws.Range("L72").Value = 33816.7503
$ws.Range("M72").Value = -5681517.540000001
$ws.Range("N72").Value = -42552.7503

$ws.Range("H129").Value = 36342
$ws.Range("I129").Value = 225844
$ws.Range("J129").Value = 6021.68
$ws.Range("K129").Value = 677532
$ws.Range("L129").Value = 18065.04
$ws.Range("M129").Value = -672532
$ws.Range("N129").Value = -28065.04

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2301.4
$ws.Range("I2").Value = 2233.125
$ws.Range("J2").Value = 2574.5
$ws.Range("K2").Value = 2233.125
$ws.Range("L2").Value = 2574.5
$ws.Range("M2").Value = -2120.125
$ws.Range("N2").Value = -2800.5

$ws.Range("H62").Value = 48124.5
$ws.Range("I62").Value = 9749.5
$ws.Range("J62").Value = 86499.5
$ws.Range("K62").Value = 9749.5
$ws.Range("L62").Value = 86499.5
$ws.Range("M62").Value = -9125.5
$ws.Range("N62").Value = -87747.5

$ws.Range("H65").Value = 48124.5
$ws.Range("I65").Value = 9749.5
$ws.Range("J65").Value = 86499.5
$ws.Range("K65").Value = 29248.5
$ws.Range("L65").Value = 259498.5
$ws.Range("M65").Value = -26128.5
$ws.Range("N65").Value = -265738.5

$ws.Range("H74").Value = 1658.9445
$ws.Range("I74").Value = 1636.5143
$ws.Range("J74").Value = 2444
$ws.Range("K74").Value = 1636.5143
$ws.Range("L74").Value = 2444
$ws.Range("M74").Value = -762.5143
$ws.Range("N74").Value = -4192

$ws.Range("H77").Value = 1658.9445
$ws.Range("I77").Value = 1636.5143
$ws.Range("J77").Value = 2444
$ws.Range("K77").Value = 8182.5715
$ws.Range("L77").Value = 12220
$ws.Range("M77").Value = -3814.5715
$ws.Range("N77").Value = -20956

$ws.Range("H116").Value = 2301.4
$ws.Range("I116").Value = 2233.125
$ws.Range("J116").Value = 2574.5
$ws.Range("K116").Value = 2233.125
$ws.Range("L116").Value = 2574.5
$ws.Range("M116").Value = 60.875
$ws.Range("N116").Value = -7162.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2301.4
$ws.Range("I3").Value = 2233.125
$ws.Range("J3").Value = 2574.5
$ws.Range("K3").Value = 2233.125
$ws.Range("L3").Value = 2574.5
$ws.Range("M3").Value = -2119.125
$ws.Range("N3").Value = -2802.5

$ws.Range("H20").Value = 3048.7778
$ws.Range("I20").Value = 2900.8333
$ws.Range("K20").Value = 2900.8333
$ws.Range("M20").Value = -2653.8333

$ws.Range("H38").Value = 27507.25
$ws.Range("J38").Value = 27507.25
$ws.Range("L38").Value = 27507.25
$ws.Range("N38").Value = -28339.25

$ws.Range("H99").Value = 2556.2856
$ws.Range("I99").Value = 1864.1177
$ws.Range("K99").Value = 1864.1177
$ws.Range("M99").Value = -366.1177

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1154.0834
$ws.Range("I16").Value = 783.6
$ws.Range("J16").Value = 3006.5
$ws.Range("K16").Value = 783.6
$ws.Range("L16").Value = 3006.5
$ws.Range("M16").Value = -496.6
$ws.Range("N16").Value = -3580.5

$ws.Range("H31").Value = 2181.2068
$ws.Range("I31").Value = 1894.6136
$ws.Range("J31").Value = 3081.9285
$ws.Range("K31").Value = 1894.6136
$ws.Range("L31").Value = 3081.9285
$ws.Range("M31").Value = -1599.6136
$ws.Range("N31").Value = -3671.9285

$ws.Range("H34").Value = 2181.2068
$ws.Range("I34").Value = 1894.6136
$ws.Range("J34").Value = 3081.9285
$ws.Range("K34").Value = 1894.6136
$ws.Range("L34").Value = 3081.9285
$ws.Range("M34").Value = -1692.6136
$ws.Range("N34").Value = -3485.9285

$ws.Range("H88").Value = 47497.25
$ws.Range("J88").Value = 47497.25
$ws.Range("L88").Value = 47497.25
$ws.Range("N88").Value = -48309.25

$ws.Range("H91").Value = 47497.25
$ws.Range("J91").Value = 47497.25
$ws.Range("L91").Value = 47497.25
$ws.Range("N91").Value = -50305.25

$ws.Range("H99").Value = 5428.4614
$ws.Range("I99").Value = 4452.1113
$ws.Range("J99").Value = 7625.25
$ws.Range("K99").Value = 4452.1113
$ws.Range("L99").Value = 7625.25
$ws.Range("M99").Value = -2954.1113
$ws.Range("N99").Value = -10621.25

$ws.Range("H109").Value = 49999
$ws.Range("J109").Value = 49999
$ws.Range("L109").Value = 49999
$ws.Range("N109").Value = -52079

$ws.Range("H113").Value = 1154.0834
$ws.Range("I113").Value = 783.6
$ws.Range("J113").Value = 3006.5
$ws.Range("K113").Value = 783.6
$ws.Range("L113").Value = 3006.5
$ws.Range("M113").Value = 1386.4
$ws.Range("N113").Value = -7346.5

$ws.Range("H122").Value = 2268.1667
$ws.Range("I122").Value = 1615.8572
$ws.Range("J122").Value = 3181.4
$ws.Range("K122").Value = 4847.571599999999
$ws.Range("L122").Value = 9544.2
$ws.Range("M122").Value = -2397.571599999999
$ws.Range("N122").Value = -14444.2

$ws.Range("H126").Value = 5428.4614
$ws.Range("I126").Value = 4452.1113
$ws.Range("J126").Value = 7625.25
$ws.Range("K126").Value = 13356.3339
$ws.Range("L126").Value = 22875.75
$ws.Range("M126").Value = -10886.3339
$ws.Range("N126").Value = -27815.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1714.4615
$ws.Range("I12").Value = 392
$ws.Range("K12").Value = 1176
$ws.Range("M12").Value = -1003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 78.708336
$ws.Range("I2").Value = 72.666664
$ws.Range("J2").Value = 88.77778
$ws.Range("K2").Value = 72.666664
$ws.Range("L2").Value = 88.77778
$ws.Range("M2").Value = 40.333336
$ws.Range("N2").Value = -314.77778

$ws.Range("H57").Value = 21127.285
$ws.Range("I57").Value = 3976
$ws.Range("K57").Value = 3976
$ws.Range("M57").Value = -3156

$ws.Range("H63").Value = 17998.8
$ws.Range("J63").Value = 18248.75
$ws.Range("L63").Value = 18248.75
$ws.Range("N63").Value = -19620.75

$ws.Range("H66").Value = 17998.8
$ws.Range("J66").Value = 18248.75
$ws.Range("L66").Value = 54746.25
$ws.Range("N66").Value = -61610.25

$ws.Range("H113").Value = 2384.4285
$ws.Range("I113").Value = 2384.4285
$ws.Range("K113").Value = 2384.4285
$ws.Range("M113").Value = -214.4285

$ws.Range("H125").Value = 67021
$ws.Range("J125").Value = 67021
$ws.Range("L125").Value = 67021
$ws.Range("N125").Value = -71941

$ws.Range("H137").Value = 92096.71
$ws.Range("J137").Value = 92096.71
$ws.Range("L137").Value = 92096.71
$ws.Range("N137").Value = -102296.71

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6342.9854
$ws.Range("I22").Value = 1699.5385
$ws.Range("J22").Value = 7440.5273
$ws.Range("K22").Value = 1699.5385
$ws.Range("L22").Value = 7440.5273
$ws.Range("M22").Value = -1404.5385
$ws.Range("N22").Value = -8030.5273

$ws.Range("H23").Value = 10794.8
$ws.Range("I23").Value = 10794.8
$ws.Range("K23").Value = 10794.8
$ws.Range("M23").Value = -10564.8

$ws.Range("H27").Value = 6342.9854
$ws.Range("I27").Value = 1699.5385
$ws.Range("J27").Value = 7440.5273
$ws.Range("K27").Value = 1699.5385
$ws.Range("L27").Value = 7440.5273
$ws.Range("M27").Value = -1592.5385
$ws.Range("N27").Value = -7654.5273

$ws.Range("H62").Value = 59969.8
$ws.Range("J62").Value = 59974.75
$ws.Range("L62").Value = 59974.75
$ws.Range("N62").Value = -61222.75

$ws.Range("H65").Value = 59969.8
$ws.Range("J65").Value = 59974.75
$ws.Range("L65").Value = 179924.25
$ws.Range("N65").Value = -186164.25

$ws.Range("H82").Value = 55557104
$ws.Range("I82").Value = 1967.6666
$ws.Range("J82").Value = 111112240
$ws.Range("K82").Value = 1967.6666
$ws.Range("L82").Value = 111112240
$ws.Range("M82").Value = -1606.6666
$ws.Range("N82").Value = -111112962

$ws.Range("H85").Value = 55557104
$ws.Range("I85").Value = 1967.6666
$ws.Range("J85").Value = 111112240
$ws.Range("K85").Value = 1967.6666
$ws.Range("L85").Value = 111112240
$ws.Range("M85").Value = -719.6666
$ws.Range("N85").Value = -111114736

$ws.Range("H132").Value = 3245.7446
$ws.Range("I132").Value = 2921.9268
$ws.Range("J132").Value = 5458.5
$ws.Range("K132").Value = 8765.7804
$ws.Range("L132").Value = 16375.5
$ws.Range("M132").Value = -6235.7804
$ws.Range("N132").Value = -21435.5

$ws.Range("H133").Value = 86569.86
$ws.Range("J133").Value = 86569.86
$ws.Range("L133").Value = 86569.86
$ws.Range("N133").Value = -91629.86

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 56499
$ws.Range("I27").Value = 26999
$ws.Range("J27").Value = 85999
$ws.Range("K27").Value = 26999
$ws.Range("L27").Value = 85999
$ws.Range("M27").Value = -26930
$ws.Range("N27").Value = -86137

$ws.Range("H122").Value = 3850.8667
$ws.Range("I122").Value = 3832.35
$ws.Range("K122").Value = 11497.05
$ws.Range("M122").Value = -9047.05

$ws.Range("H126").Value = 20267.955
$ws.Range("I126").Value = 20267.955
$ws.Range("K126").Value = 60803.86500000001
$ws.Range("M126").Value = -58333.86500000001

$ws.Range("H132").Value = 4514.2
$ws.Range("I132").Value = 3121.05
$ws.Range("J132").Value = 7300.5
$ws.Range("K132").Value = 9363.150000000001
$ws.Range("L132").Value = 21901.5
$ws.Range("M132").Value = -6833.150000000001
$ws.Range("N132").Value = -26961.5

$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120
